# Open issues sheet: update statuses and owner, fix selection.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Foglio1")

# Status column (D) updates
$ws.Range("D2").Value = "Chiuso"
$ws.Range("D3").Value = "Chiuso"
$ws.Range("D5").Value = "In corso"
$ws.Range("D6").Value = "Chiuso"
$ws.Range("D15").Value = "Chiuso"

# Owner column (F) update for row 5
$ws.Range("F5").Value = "Marco / Daniele"

# The longer owner text no longer fits column F's default width; re-apply
# best-fit auto sizing like Excel does automatically on entry.
$ws.Columns.Item(6).AutoFit()

# Restore the active selection to D14 as in the author's edit
$ws.Activate()
$ws.Range("D14").Select()
